$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-01-08 Wednesday" "2025-01-14 Tuesday"

Replace-Text "654÷7=93, 3" "581÷6=96, 5"
Replace-Text "832÷9=92, 4" "248÷7=35, 3"
Replace-Text "339÷6=56, 3" "365÷2=182, 1"
Replace-Text "131÷3=43, 2" "749÷9=83, 2"
Replace-Text "333÷4=83, 1" "128÷8=16, 0"

Replace-Text "486÷6=81, 0" "680÷6=113, 2"
Replace-Text "343÷2=171, 1" "543÷3=181, 0"
Replace-Text "246÷4=61, 2" "300÷8=37, 4"
Replace-Text "696÷2=348, 0" "467÷6=77, 5"
Replace-Text "210÷3=70, 0" "728÷6=121, 2"

Replace-Text "332÷7=47, 3" "453÷4=113, 1"
Replace-Text "956÷4=239, 0" "234÷7=33, 3"
Replace-Text "783÷7=111, 6" "544÷8=68, 0"
Replace-Text "410÷4=102, 2" "626÷4=156, 2"
Replace-Text "869÷8=108, 5" "345÷5=69, 0"

Replace-Text "777÷6=129, 3" "667÷7=95, 2"
Replace-Text "302÷8=37, 6" "912÷4=228, 0"
Replace-Text "645÷2=322, 1" "919÷9=102, 1"
Replace-Text "938÷3=312, 2" "548÷5=109, 3"
Replace-Text "411÷6=68, 3" "310÷8=38, 6"

Replace-Text "609÷8=76, 1" "920÷4=230, 0"
Replace-Text "487÷6=81, 1" "570÷8=71, 2"
Replace-Text "273÷2=136, 1" "458÷2=229, 0"
Replace-Text "499÷2=249, 1" "195÷4=48, 3"
Replace-Text "519÷2=259, 1" "748÷7=106, 6"

Write-Output "done"
